$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.03958866666666667
$ws.Range("H2").Value = 0.118766
$ws.Range("I2").Value = 0.0007442768123675562
$ws.Range("J2").Value = 0.0007442768123675561
$ws.Range("M2").Value = 4.717738333333333
$ws.Range("N2").Value = 14.153215
$ws.Range("O2").Value = 0.2002263444295212
$ws.Range("P2").Value = 0.2002263444295212
$ws.Range("Q2").Value = 0.1867689702988889
$ws.Range("R2").Value = 1.68092073269
$ws.Range("S2").Value = 0.0001490238253840124
$ws.Range("T2").Value = 0.0001490238253840124

$ws.Range("G3").Value = 0.03958866666666667
$ws.Range("H3").Value = 0.118766
$ws.Range("I3").Value = 0.0007442768123675562
$ws.Range("J3").Value = 0.0007442768123675561
$ws.Range("N3").Value = 9.228847
$ws.Range("O3").Value = 0.1305610278731266
$ws.Range("P3").Value = 0.1305610278731266
$ws.Range("Q3").Value = 0.1217859158668889
$ws.Range("R3").Value = 1.096073242802
$ws.Range("S3").Value = 0.00009717354564484232
$ws.Range("T3").Value = 0.00009717354564484231

$ws.Range("G4").Value = 0.03958866666666667
$ws.Range("H4").Value = 0.118766
$ws.Range("I4").Value = 0.0007442768123675562
$ws.Range("J4").Value = 0.0007442768123675561
$ws.Range("M4").Value = 0.6908423333333333
$ws.Range("N4").Value = 2.072527
$ws.Range("O4").Value = 0.02932015834857891
$ws.Range("P4").Value = 0.02932015834857891
$ws.Range("Q4").Value = 0.02734952685355556
$ws.Range("R4").Value = 0.246145741682
$ws.Range("S4").Value = 0.0000218223139937923
$ws.Range("T4").Value = 0.0000218223139937923

$ws.Range("G5").Value = 0.03958866666666667
$ws.Range("H5").Value = 0.118766
$ws.Range("I5").Value = 0.0007442768123675562
$ws.Range("J5").Value = 0.0007442768123675561
$ws.Range("M5").Value = 15.077163
$ws.Range("N5").Value = 45.231489
$ws.Range("O5").Value = 0.6398924693487733
$ws.Range("P5").Value = 0.6398924693487733
$ws.Range("Q5").Value = 0.596884780286
$ws.Range("R5").Value = 5.371963022574
$ws.Range("S5").Value = 0.0004762571273449092
$ws.Range("T5").Value = 0.0004762571273449091

$ws.Range("H6").Value = 0.059669
$ws.Range("I6").Value = 0.0003739306966401135
$ws.Range("J6").Value = 0.0003739306966401134
$ws.Range("M6").Value = 4.717738333333333
$ws.Range("N6").Value = 14.153215
$ws.Range("O6").Value = 0.2002263444295212
$ws.Range("P6").Value = 0.2002263444295212
$ws.Range("Q6").Value = 0.09383424287055556
$ws.Range("R6").Value = 0.844508185835
$ws.Range("S6").Value = 0.00007487077645823415
$ws.Range("T6").Value = 0.00007487077645823415

$ws.Range("H7").Value = 0.059669
$ws.Range("I7").Value = 0.0003739306966401135
$ws.Range("J7").Value = 0.0003739306966401134
$ws.Range("N7").Value = 9.228847
$ws.Range("O7").Value = 0.1305610278731266
$ws.Range("P7").Value = 0.1305610278731266
$ws.Range("R7").Value = 0.550676071643
$ws.Range("S7").Value = 0.0000488207761066475
$ws.Range("T7").Value = 0.00004882077610664749

$ws.Range("H8").Value = 0.059669
$ws.Range("I8").Value = 0.0003739306966401135
$ws.Range("J8").Value = 0.0003739306966401134
$ws.Range("M8").Value = 0.6908423333333333
$ws.Range("N8").Value = 2.072527
$ws.Range("O8").Value = 0.02932015834857891
$ws.Range("P8").Value = 0.02932015834857891
$ws.Range("Q8").Value = 0.01374062372922222
$ws.Range("R8").Value = 0.123665613563
$ws.Range("S8").Value = 0.00001096370723688255
$ws.Range("T8").Value = 0.00001096370723688255

$ws.Range("H9").Value = 0.059669
$ws.Range("I9").Value = 0.0003739306966401135
$ws.Range("J9").Value = 0.0003739306966401134
$ws.Range("M9").Value = 15.077163
$ws.Range("N9").Value = 45.231489
$ws.Range("O9").Value = 0.6398924693487733
$ws.Range("P9").Value = 0.6398924693487733
$ws.Range("Q9").Value = 0.299879746349
$ws.Range("R9").Value = 2.698917717141
$ws.Range("S9").Value = 0.0002392754368383492
$ws.Range("T9").Value = 0.0002392754368383492

$ws.Range("G10").Value = 53.131305
$ws.Range("H10").Value = 159.393915
$ws.Range("I10").Value = 0.9988817924909924
$ws.Range("J10").Value = 0.9988817924909923
$ws.Range("M10").Value = 4.717738333333333
$ws.Range("N10").Value = 14.153215
$ws.Range("O10").Value = 0.2002263444295212
$ws.Range("P10").Value = 0.2002263444295212
$ws.Range("Q10").Value = 250.659594298525
$ws.Range("R10").Value = 2255.936348686725
$ws.Range("S10").Value = 0.2000024498276789
$ws.Range("T10").Value = 0.200002449827679

$ws.Range("G11").Value = 53.131305
$ws.Range("H11").Value = 159.393915
$ws.Range("I11").Value = 0.9988817924909924
$ws.Range("J11").Value = 0.9988817924909923
$ws.Range("N11").Value = 9.228847
$ws.Range("O11").Value = 0.1305610278731266
$ws.Range("P11").Value = 0.1305610278731266
$ws.Range("Q11").Value = 163.446894918445
$ws.Range("R11").Value = 1471.022054266005
$ws.Range("S11").Value = 0.1304150335513751
$ws.Range("T11").Value = 0.1304150335513751

$ws.Range("G12").Value = 53.131305
$ws.Range("H12").Value = 159.393915
$ws.Range("I12").Value = 0.9988817924909924
$ws.Range("J12").Value = 0.9988817924909923
$ws.Range("M12").Value = 0.6908423333333333
$ws.Range("N12").Value = 2.072527
$ws.Range("O12").Value = 0.02932015834857891
$ws.Range("P12").Value = 0.02932015834857891
$ws.Range("Q12").Value = 36.705354719245
$ws.Range("R12").Value = 330.348192473205
$ws.Range("S12").Value = 0.02928737232734824
$ws.Range("T12").Value = 0.02928737232734823

$ws.Range("G13").Value = 53.131305
$ws.Range("H13").Value = 159.393915
$ws.Range("I13").Value = 0.9988817924909924
$ws.Range("J13").Value = 0.9988817924909923
$ws.Range("M13").Value = 15.077163
$ws.Range("N13").Value = 45.231489
$ws.Range("O13").Value = 0.6398924693487733
$ws.Range("P13").Value = 0.6398924693487733
$ws.Range("Q13").Value = 801.0693458877149
$ws.Range("R13").Value = 7209.624112989434
$ws.Range("S13").Value = 0.6391769367845901
$ws.Range("T13").Value = 0.63917693678459
